$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.377.92'
$ws.Range('E2').Value = '  -0.37%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.566.54'
$ws.Range('E3').Value = '  -0.08%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.97'
$ws.Range('E5').Value = '  -0.46%  '

$ws.Range('E6').Value = '  -0.66%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.39'
$ws.Range('E8').Value = '  -3.53%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.51'
$ws.Range('E9').Value = '  -2.12%  '

$ws.Range('E10').Value = '  -1.41%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0588'
$ws.Range('E11').Value = '  -0.71%  '

$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.788.41'
$ws.Range('E13').Value = '  -0.20%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.572.26'
$ws.Range('E14').Value = '  +0.32%  '

$ws.Range('E15').Value = '  -0.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.362.06'
$ws.Range('E16').Value = '  -0.48%  '

$ws.Range('E17').Value = '  -1.62%  '

$ws.Range('E18').Value = '  -3.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.95'

$ws.Range('E20').Value = '  +0.33%  '

$ws.Range('E21').Value = '  -1.95%  '

$ws.Range('E22').Value = '  +0.04%  '

$ws.Range('E24').Value = '  -2.26%  '

$ws.Range('E25').Value = '  -2.15%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.06'
$ws.Range('E26').Value = '  -0.64%  '

$ws.Range('E27').Value = '  -0.85%  '

$ws.Range('E28').Value = '  +0.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.31'
$ws.Range('E29').Value = '  -2.32%  '

$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('E31').Value = '  +1.76%  '

$ws.Range('E32').Value = '  -3.76%  '

$ws.Range('E33').Value = '  -1.14%  '

$ws.Range('E34').Value = '  -0.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.383.66'
$ws.Range('E35').Value = '  -0.96%  '

$ws.Range('E36').Value = '  +1.80%  '

$ws.Range('E37').Value = '  -3.31%  '

$ws.Range('E38').Value = '  -0.56%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.65'
$ws.Range('E39').Value = '  +2.07%  '

$ws.Range('E40').Value = '  -2.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.95'
$ws.Range('E41').Value = '  +3.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.518'
$ws.Range('E42').Value = '  -3.08%  '

$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0467'
$ws.Range('E45').Value = '  -1.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.34'
$ws.Range('E46').Value = '  -2.94%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.922'
$ws.Range('E47').Value = '  -5.31%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '62.20'
$ws.Range('E48').Value = '  -1.03%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.701.47'
$ws.Range('E49').Value = '  -0.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.34'
$ws.Range('E50').Value = '  -0.70%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0515'
$ws.Range('E51').Value = '  -2.18%  '
